$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely - remaining rows shift up by one.
$ws.Rows("26:26").Delete()

# After the above deletion, the row that used to be "SC 92" (old row 28) is
# now at row 27. Delete it entirely - remaining rows shift up by one more.
$ws.Rows("27:27").Delete()

# Now rows 26-33 hold (in order): SC 5, SC 101, SC 105, SC 119, SC 120,
# SC 132, SC 193, SC 232. Update column C ("B" header) values to match the
# target data - some previously-missing values are now filled in, and some
# previously-filled values are now blanked out.
$ws.Range("C26").Value = 10.8
$ws.Range("C27").Value = ""
$ws.Range("C28").Value = ""
$ws.Range("C29").Value = 11.2
$ws.Range("C30").Value = 11.4
$ws.Range("C31").Value = ""
$ws.Range("C32").Value = ""
$ws.Range("C33").Value = 10.4
